# Update "想去人数" (want-to-go count) values in the "展览" and "全部类型" sheets
# to reflect newly generated output numbers.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 6531
$ws1.Range("F3").Value = 120
$ws1.Range("F5").Value = 405
$ws1.Range("F9").Value = 94
$ws1.Range("F10").Value = 86
$ws1.Range("F11").Value = 129
$ws1.Range("F13").Value = 384
$ws1.Range("F14").Value = 1112
$ws1.Range("F15").Value = 3231
$ws1.Range("F17").Value = 202
$ws1.Range("F18").Value = 1884

# Sheet "全部类型" (all types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 6531
$ws4.Range("F3").Value = 120
$ws4.Range("F5").Value = 405
$ws4.Range("F10").Value = 94
$ws4.Range("F11").Value = 86
$ws4.Range("F12").Value = 129
$ws4.Range("F14").Value = 384
$ws4.Range("F15").Value = 1112
$ws4.Range("F16").Value = 3231
$ws4.Range("F18").Value = 202
$ws4.Range("F19").Value = 1884
